# Daily attendance processing - 2026-01-06 17:09:58
#
# Normalizes the "Recorded By" (column G) cell text on the
# "Session Analysis Results" sheet: for the specific recorder-name
# combinations produced by the nightly export, the first name in the
# comma-separated list is moved to the end of the list (a left-rotate
# of the list by one position). Combinations not in the lookup table
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact-match lookup table: old "Recorded By" text -> new "Recorded By" text.
$rotations = @{
    "dnasr281@gmail.com, System"                       = "System, dnasr281@gmail.com";
    "system, System, backup@backdoor.com"               = "System, backup@backdoor.com, system";
    "System, admin@admin.com"                           = "admin@admin.com, System";
    "dnasr281@gmail.com, admin@admin.com"               = "admin@admin.com, dnasr281@gmail.com";
}

$lastRow = $ws.Cells.Item(1, 1).SpecialCells(11).Row
$col = 7  # column G = "Recorded By"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $current = $cell.Text

    if ($rotations.ContainsKey($current)) {
        $cell.Value = $rotations[$current]
    }
}
